$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 12:08"

# --- Row 37 (Rumania): updated totals, no name change ---
$ws.Range("B37").Value = 13512
$ws.Range("C37").Value = 349
$ws.Range("D37").Value = 5269
$ws.Range("E37").Value = 7440
$ws.Range("F37").Value = 243

# --- Row 59 (Moldavia): updated totals, no name change ---
$ws.Range("D59").Value = 1423
$ws.Range("E59").Value = 2570
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 128

# --- Row 97 (Albania): updated totals, no name change ---
$ws.Range("B97").Value = 803
$ws.Range("C97").Value = 8
$ws.Range("D97").Value = 543
$ws.Range("E97").Value = 229

# --- Rows 144/145: Brunei / Etiopia swap places (Etiopia overtakes Brunei) ---
$ws.Range("A144").Value = "Etiopia"
$ws.Range("B144").Value = 140
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 75
$ws.Range("E144").Value = 62
$ws.Range("F144").Value = 0
$ws.Range("H144").Value = 3

$ws.Range("A145").Value = "Brunei"
$ws.Range("B145").Value = 138
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 130
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = 2
$ws.Range("H145").Value = 1

# --- Rows 188/189: Belice / Santa Lucia swap places ---
$ws.Range("A188").Value = "Santa Lucia"
$ws.Range("D188").Value = 15
$ws.Range("F188").Value = 0
$ws.Range("H188").Value = 0

$ws.Range("A189").Value = "Belice"
$ws.Range("D189").Value = 13
$ws.Range("F189").Value = 1
$ws.Range("H189").Value = 2

# --- Rows 194/195: San Vicente y las Granadinas / Namibia swap places (values identical) ---
$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

# --- Rows 198/199: San Cristobal y Nieves / Burundi swap places ---
$ws.Range("A198").Value = "Burundi"
$ws.Range("D198").Value = 7
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0
